$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for the "Out of PO" roster table (player, position, team),
# replacing Harrison Barnes (San Antonio Spurs) with Devin Vassell
# (San Antonio Spurs) and re-ordering the remaining rows.
$data = @(
    @("Toumani Camara",     "SG,SF,PF", "Portland Trail Blazers"),
    @("Kelly Oubre Jr.",    "SG,SF",    "Philadelphia 76ers"),
    @("Michael Porter Jr.", "SF,PF",    "Denver Nuggets"),
    @("Alperen Sengün",     "C",        "Houston Rockets"),
    @("Gary Trent Jr.",     "PG,SG,SF", "Milwaukee Bucks"),
    @("Domantas Sabonis",   "C",        "Sacramento Kings"),
    @("Max Strus",          "SG,SF",    "Cleveland Cavaliers"),
    @("Malik Beasley",      "SG,SF",    "Detroit Pistons"),
    @("Donovan Mitchell",   "PG,SG",    "Cleveland Cavaliers"),
    @("Bilal Coulibaly",    "SG,SF",    "Washington Wizards"),
    @("Jaden McDaniels",    "SF,PF",    "Minnesota Timberwolves"),
    @("De'Andre Hunter",    "SF,PF",    "Cleveland Cavaliers"),
    @("Kristaps Porzingis", "PF,C",     "Boston Celtics"),
    @("Dyson Daniels",      "PG,SG,SF", "Atlanta Hawks"),
    @("Devin Vassell",      "SG,SF",    "San Antonio Spurs"),
    @("Julius Randle",      "PF,C",     "Minnesota Timberwolves"),
    @("Cam Thomas",         "SG,SF",    "Brooklyn Nets"),
    @("Josh Hart",          "SG,SF,PF", "New York Knicks")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
